$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it a
# genuine text cell (Coinranking's "Price" column stores numbers as text,
# e.g. "410.73"). Plain `.Value = "410.73"` lets Excel auto-coerce the
# string to a real number, so we briefly force Text format, assign, then
# drop the format again (ClearFormats) to avoid leaving a stray style
# behind while still keeping the cell's stored type as text.
function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "62.027.87"
$ws.Range("D3").Value = "3.417.39"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue $ws "D5" "410.73"
Set-TextValue $ws "D6" "129.35"
$ws.Range("E6").Value = "  -3.93%  "
Set-TextValue $ws "D7" "0.643"
$ws.Range("E7").Value = "  +8.72%  "
Set-TextValue $ws "D8" "1.00"
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue $ws "D9" "0.739"
$ws.Range("E9").Value = "  +7.76%  "
Set-TextValue $ws "D10" "0.141"
$ws.Range("E10").Value = "  +15.46%  "
Set-TextValue $ws "D11" "42.77"
$ws.Range("E11").Value = "  +0.98%  "
Set-TextValue $ws "D12" "0.0000217"
$ws.Range("E12").Value = "  +65.36%  "
Set-TextValue $ws "D13" "9.15"
$ws.Range("E13").Value = "  +8.05%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "3.953.63"
$ws.Range("E15").Value = "  -0.56%  "
Set-TextValue $ws "D16" "21.20"
$ws.Range("E16").Value = "  +6.39%  "
$ws.Range("D17").Value = "3.401.57"
$ws.Range("E17").Value = "  -1.24%  "
Set-TextValue $ws "D18" "12.07"
$ws.Range("E18").Value = "  +5.93%  "
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("D20").Value = "61.987.83"
$ws.Range("E20").Value = "  -0.49%  "
Set-TextValue $ws "D21" "445.48"
$ws.Range("E21").Value = "  +41.46%  "
Set-TextValue $ws "D22" "91.42"
$ws.Range("E22").Value = "  +8.53%  "
$ws.Range("E23").Value = "  -0.50%  "
Set-TextValue $ws "D24" "13.13"
$ws.Range("E24").Value = "  +1.32%  "
Set-TextValue $ws "D25" "3.26"
$ws.Range("E25").Value = "  +3.11%  "
Set-TextValue $ws "D26" "33.49"
$ws.Range("E26").Value = "  +12.44%  "
Set-TextValue $ws "D27" "8.84"
$ws.Range("E27").Value = "  +7.74%  "
$ws.Range("E28").Value = "  +0.06%  "
Set-TextValue $ws "D29" "7.61"
$ws.Range("E29").Value = "  +0.41%  "
Set-TextValue $ws "D30" "2.74"
$ws.Range("E30").Value = "  +0.15%  "
Set-TextValue $ws "D31" "11.98"
$ws.Range("E31").Value = "  +5.35%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  -2.45%  "
Set-TextValue $ws "D34" "42.92"
Set-TextValue $ws "D35" "0.999"
$ws.Range("E35").Value = "  -0.09%  "
Set-TextValue $ws "D36" "0.0502"
$ws.Range("E36").Value = "  +3.61%  "
Set-TextValue $ws "D37" "53.80"
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  +8.48%  "
Set-TextValue $ws "D40" "3.37"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -1.47%  "
Set-TextValue $ws "D42" "0.315"
$ws.Range("E42").Value = "  +0.63%  "
Set-TextValue $ws "D43" "141.13"
$ws.Range("E43").Value = "  +2.29%  "
Set-TextValue $ws "D44" "4.23"
$ws.Range("E44").Value = "  +4.86%  "
Set-TextValue $ws "D45" "1.98"
$ws.Range("E45").Value = "  -0.35%  "
Set-TextValue $ws "D46" "2.41"
$ws.Range("E46").Value = "  +7.63%  "
Set-TextValue $ws "D47" "16.63"
$ws.Range("E47").Value = "  -1.03%  "
Set-TextValue $ws "D48" "22.19"
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").Value = "3.766.90"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "2.115.18"
$ws.Range("E50").Value = "  -0.56%  "
Set-TextValue $ws "D51" "104.87"
$ws.Range("E51").Value = "  +25.21%  "